$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E width change (22 -> 20)
# Note: the host's ColumnWidth setter adds a constant +5/6 offset when
# round-tripping through the OOXML <col width> attribute, so we
# pre-compensate to land exactly on 20.
$ws.Columns.Item(5).ColumnWidth = 19.166666666666668

# Row 2
$ws.Range("C2").Value = 43460
$ws.Range("D2").Value = "Visa"
$ws.Range("E2").Value = "Walgreens"
$ws.Range("F2").Value = "Pharmacy"
$ws.Range("G2").Value = "General"
$ws.Range("I2").Value = -359.77
$ws.Range("J2").Value = "Expense"
$ws.Range("K2").Value = "Discretionary"

# Row 3
$ws.Range("C3").Value = 43471
$ws.Range("E3").Value = "Amazon"
$ws.Range("F3").Value = "Shopping"
$ws.Range("G3").Value = "Online"
$ws.Range("I3").Value = -430.68

# Row 4
$ws.Range("C4").Value = 43569
$ws.Range("D4").Value = "Mastercard"
$ws.Range("E4").Value = "Neighborhood Market"
$ws.Range("F4").Value = "Groceries"
$ws.Range("G4").Value = "Farmer's Market"
$ws.Range("I4").Value = -417.02
$ws.Range("K4").Value = "Essential"

# Row 5
$ws.Range("C5").Value = 43726
$ws.Range("D5").Value = "Mastercard"
$ws.Range("E5").Value = "Kroger"
$ws.Range("F5").Value = "Groceries"
$ws.Range("G5").Value = "Grocery Store"
$ws.Range("I5").Value = -423.43
$ws.Range("J5").Value = "Expense"
$ws.Range("K5").Value = "Essential"

# Row 6
$ws.Range("C6").Value = 43987
$ws.Range("D6").Value = "Visa"
$ws.Range("E6").Value = "Neighborhood Market"
$ws.Range("G6").Value = "Farmer's Market"
$ws.Range("I6").Value = -157

# Row 7
$ws.Range("C7").Value = 44094
$ws.Range("D7").Value = "Savings"
$ws.Range("E7").Value = "Wells Fargo"
$ws.Range("F7").Value = "Interest Income"
$ws.Range("G7").Value = ""
$ws.Range("I7").Value = 3608.63
$ws.Range("J7").Value = "Income"
$ws.Range("K7").Value = "Income"

# Row 8
$ws.Range("C8").Value = 44223
$ws.Range("D8").Value = "Mastercard"
$ws.Range("E8").Value = "McDonald's"
$ws.Range("F8").Value = "Dining Out"
$ws.Range("I8").Value = -159.67
$ws.Range("J8").Value = "Expense"
$ws.Range("K8").Value = "Discretionary"
$ws.Range("O8").Value = "Era A"

# Row 9
$ws.Range("C9").Value = 44226
$ws.Range("D9").Value = "Savings"
$ws.Range("E9").Value = "Wells Fargo"
$ws.Range("F9").Value = "Interest Income"
$ws.Range("G9").Value = ""
$ws.Range("I9").Value = 3849.4
$ws.Range("J9").Value = "Income"
$ws.Range("K9").Value = "Income"
$ws.Range("O9").Value = "Era A"

# Row 10
$ws.Range("C10").Value = 44517
$ws.Range("D10").Value = "Mastercard"
$ws.Range("E10").Value = "The Vet"
$ws.Range("F10").Value = "Pets"
$ws.Range("G10").Value = ""
$ws.Range("I10").Value = -16.17

# Row 11
$ws.Range("C11").Value = 44799
$ws.Range("E11").Value = "Highland Apartments"
$ws.Range("F11").Value = "Rent"
$ws.Range("I11").Value = -363.44
$ws.Range("K11").Value = "Discretionary"

# Row 12
$ws.Range("C12").Value = 44879
$ws.Range("D12").Value = "Savings"
$ws.Range("E12").Value = "Wells Fargo"
$ws.Range("F12").Value = "Interest Income"
$ws.Range("G12").Value = ""
$ws.Range("I12").Value = 1949.29
$ws.Range("J12").Value = "Income"
$ws.Range("K12").Value = "Income"
